$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.957.51'
$ws.Range('E2').Value = '  -0.83%  '

$ws.Range('D3').Value = '3.531.61'
$ws.Range('E3').Value = '  -1.27%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').Value = '608.38'
$ws.Range('E5').Value = '  +2.96%  '

$ws.Range('D6').Value = '184.85'
$ws.Range('E6').Value = '  -1.69%  '

$ws.Range('D7').Value = '3.526.71'
$ws.Range('E7').Value = '  -1.18%  '

$ws.Range('D8').Value = '0.612'
$ws.Range('E8').Value = '  -1.64%  '

$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.03%  '

$ws.Range('D10').Value = '0.212'
$ws.Range('E10').Value = '  +6.17%  '

$ws.Range('D11').Value = '0.640'
$ws.Range('E11').Value = '  -1.49%  '

$ws.Range('D12').Value = '53.47'
$ws.Range('E12').Value = '  -2.78%  '

$ws.Range('D13').Value = '0.0000307'
$ws.Range('E13').Value = '  -0.96%  '

$ws.Range('D14').Value = '9.41'
$ws.Range('E14').Value = '  -2.17%  '

$ws.Range('D15').Value = '4.083.01'
$ws.Range('E15').Value = '  -1.39%  '

$ws.Range('D16').Value = '69.955.72'
$ws.Range('E16').Value = '  -0.65%  '

$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = '18.86'
$ws.Range('E17').Value = '  -3.37%  '

$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = '12.60'
$ws.Range('E18').Value = '  +0.74%  '

$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.526.90'
$ws.Range('E19').Value = '  -1.12%  '

$ws.Range('D20').Value = '576.38'
$ws.Range('E20').Value = '  +2.65%  '

$ws.Range('E21').Value = '  -0.03%  '

$ws.Range('D22').Value = '0.988'
$ws.Range('E22').Value = '  -3.49%  '

$ws.Range('D23').Value = '17.29'
$ws.Range('E23').Value = '  -4.11%  '

$ws.Range('D24').Value = '4.69'
$ws.Range('E24').Value = '  +0.12%  '

$ws.Range('D25').Value = '4.83'
$ws.Range('E25').Value = '  -2.40%  '

$ws.Range('D26').Value = '93.76'
$ws.Range('E26').Value = '  -2.58%  '

$ws.Range('D27').Value = '11.04'
$ws.Range('E27').Value = '  -3.91%  '

$ws.Range('D28').Value = '2.94'
$ws.Range('E28').Value = '  -2.16%  '

$ws.Range('D29').Value = '9.30'
$ws.Range('E29').Value = '  +1.22%  '

$ws.Range('D30').Value = '32.00'
$ws.Range('E30').Value = '  -1.21%  '

$ws.Range('D31').Value = '6.99'
$ws.Range('E31').Value = '  -5.13%  '

$ws.Range('D32').Value = '12.17'
$ws.Range('E32').Value = '  -3.42%  '

$ws.Range('E33').Value = '  -1.39%  '

$ws.Range('D34').Value = '63.21'
$ws.Range('E34').Value = '  -2.99%  '

$ws.Range('D35').Value = '3.36'
$ws.Range('E35').Value = '  +2.48%  '

$ws.Range('E36').Value = '  +16.12%  '

$ws.Range('D37').Value = '531.17'
$ws.Range('E37').Value = '  -4.18%  '

$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').Value = '0.403'
$ws.Range('E38').Value = '  -3.57%  '

$ws.Range('B39').Value = 'Dai'
$ws.Range('C39').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.20%  '

$ws.Range('D40').Value = '36.97'
$ws.Range('E40').Value = '  -3.40%  '

$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0778'
$ws.Range('E41').Value = '  +0.56%  '

$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '3.530.70'
$ws.Range('E42').Value = '  +4.65%  '

$ws.Range('D43').Value = '3.51'
$ws.Range('E43').Value = '  +3.13%  '

$ws.Range('D44').Value = '0.136'
$ws.Range('E44').Value = '  +0.53%  '

$ws.Range('D45').Value = '0.0453'
$ws.Range('E45').Value = '  +1.18%  '

$ws.Range('D46').Value = '3.42'
$ws.Range('E46').Value = '  -4.88%  '

$ws.Range('D47').Value = '2.90'
$ws.Range('E47').Value = '  -3.08%  '

$ws.Range('D48').Value = '0.140'
$ws.Range('E48').Value = '  +2.56%  '

$ws.Range('D49').Value = '9.09'

$ws.Range('E50').Value = '  +0.02%  '

$ws.Range('D51').Value = '1.44'
$ws.Range('E51').Value = '  -2.26%  '
